$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue $ws "D2" "58.333.97"
Set-TextValue $ws "E2" "  +3.42%  "
Set-TextValue $ws "D3" "2.364.04"
Set-TextValue $ws "E3" "  +1.64%  "
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  -0.04%  "
Set-TextValue $ws "D5" "544.13"
Set-TextValue $ws "E5" "  +6.35%  "
Set-TextValue $ws "D6" "135.11"
Set-TextValue $ws "E6" "  +2.72%  "
Set-TextValue $ws "D7" "1.00"
Set-TextValue $ws "E7" "  +0.07%  "
Set-TextValue $ws "E8" "  +0.94%  "
Set-TextValue $ws "D9" "2.363.46"
Set-TextValue $ws "E9" "  +1.47%  "
Set-TextValue $ws "E10" "  +1.98%  "
Set-TextValue $ws "E11" "  +1.20%  "
Set-TextValue $ws "D12" "5.40"
Set-TextValue $ws "E12" "  +2.54%  "
Set-TextValue $ws "D13" "0.357"
Set-TextValue $ws "E13" "  +5.94%  "
Set-TextValue $ws "B14" "WrappedliquidstakedEther2.0"
Set-TextValue $ws "C14" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws "D14" "2.785.13"
Set-TextValue $ws "E14" "  +1.58%  "
Set-TextValue $ws "B15" "Avalanche"
Set-TextValue $ws "C15" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws "D15" "23.66"
Set-TextValue $ws "E15" "  +0.73%  "
Set-TextValue $ws "D16" "58.265.07"
Set-TextValue $ws "E16" "  +3.33%  "
Set-TextValue $ws "E17" "  +1.46%  "
Set-TextValue $ws "D18" "2.371.32"
Set-TextValue $ws "E18" "  +1.72%  "
Set-TextValue $ws "D19" "10.61"
Set-TextValue $ws "E19" "  +1.80%  "
Set-TextValue $ws "D20" "335.78"
Set-TextValue $ws "E20" "  +3.10%  "
Set-TextValue $ws "E21" "  +2.32%  "
Set-TextValue $ws "D22" "6.74"
Set-TextValue $ws "E22" "  +0.23%  "
Set-TextValue $ws "D23" "1.00"
Set-TextValue $ws "E23" "  +0.04%  "
Set-TextValue $ws "D24" "62.02"
Set-TextValue $ws "E24" "  +0.91%  "
Set-TextValue $ws "E25" "  +5.18%  "
Set-TextValue $ws "B26" "Binance-PegBSC-USD"
Set-TextValue $ws "C26" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws "D26" "1.00"
Set-TextValue $ws "E26" "  +0.15%  "
Set-TextValue $ws "B27" "InternetComputer(DFINITY)"
Set-TextValue $ws "C27" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D27" "8.48"
Set-TextValue $ws "E27" "  -3.02%  "
Set-TextValue $ws "D28" "1.43"
Set-TextValue $ws "E28" "  +10.80%  "
Set-TextValue $ws "E29" "  +5.62%  "
Set-TextValue $ws "D30" "170.58"
Set-TextValue $ws "E30" "  +1.78%  "
Set-TextValue $ws "D31" "0.0₃0739"
Set-TextValue $ws "E31" "  +3.14%  "
Set-TextValue $ws "D32" "6.16"
Set-TextValue $ws "E32" "  +1.00%  "
Set-TextValue $ws "D33" "18.57"
Set-TextValue $ws "E33" "  +1.39%  "
Set-TextValue $ws "E34" "  +15.24%  "
Set-TextValue $ws "D35" "0.999"
Set-TextValue $ws "E35" "  +0.00%  "
Set-TextValue $ws "E36" "  +0.13%  "
Set-TextValue $ws "D37" "4.20"
Set-TextValue $ws "E37" "  +7.58%  "
Set-TextValue $ws "E38" "  +1.16%  "
Set-TextValue $ws "E39" "  +6.42%  "
Set-TextValue $ws "D40" "39.36"
Set-TextValue $ws "E40" "  +2.45%  "
Set-TextValue $ws "D41" "150.59"
Set-TextValue $ws "E41" "  +0.97%  "
Set-TextValue $ws "D42" "0.380"
Set-TextValue $ws "E42" "  +1.91%  "
Set-TextValue $ws "E43" "  +2.62%  "
Set-TextValue $ws "D44" "287.08"
Set-TextValue $ws "E44" "  +3.54%  "
Set-TextValue $ws "D45" "19.34"
Set-TextValue $ws "E45" "  +7.14%  "
Set-TextValue $ws "E46" "  +0.46%  "
Set-TextValue $ws "E47" "  +3.10%  "
Set-TextValue $ws "E49" "  +2.69%  "
Set-TextValue $ws "D50" "17.70"
Set-TextValue $ws "E50" "  +3.95%  "
Set-TextValue $ws "D51" "0.381"
Set-TextValue $ws "E51" "  +1.05%  "
